# Fill in the missing meteorological data for 25-05-2025 .. 31-05-2025.
#
# The worksheet currently has a single blank "spacer" row (row 33, style s=3)
# sitting where the 25-31 May data should be, immediately followed by the
# "KETERANGAN:" notes block (rows 36-48). We need to:
#   1) push the spacer row (and the notes block below it) down by 7 rows so
#      there is room for the 7 missing days of data,
#   2) fill rows 33-39 with the real readings for 25-05-2025 .. 31-05-2025,
#      formatted like the existing data rows (thin box border, centered,
#      wrapped text - matching the other data rows' style),
#   3) leave the (now relocated) spacer row at row 40 blank, and
#   4) move the active selection to the new spacer row (A40:K40), mirroring
#      where the cursor ends up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert 7 fresh rows above the old spacer row -----------------------
# This shifts the blank spacer row from 33 down to 40, and the KETERANGAN
# block (old rows 36-48) down to rows 43-55, exactly as in the target.
$ws.Rows("33:39").Insert()

# --- 2) Style the newly inserted rows like the other data rows -------------
$dataRange = $ws.Range("A33:K39")
$dataRange.HorizontalAlignment = -4108   # xlCenter
$dataRange.VerticalAlignment = -4108     # xlCenter
$dataRange.WrapText = $true
$dataRange.Borders.LineStyle = 1         # xlContinuous
$dataRange.Borders.Weight = 2            # xlThin
$dataRange.Borders.Color = 0             # black

# --- 3) Populate the readings ----------------------------------------------
# Columns: A=TANGGAL B=TN C=TX D=TAVG E=RH_AVG F=RR G=SS H=FF_X I=DDD_X J=FF_AVG K=DDD_CAR
$data = @(
    @(33, "25-05-2025", 25.8, 33.9, 28.5, 71.0, 0.0, 8.0, 3.0, 30, 1.0, "C"),
    @(34, "26-05-2025", 24.6, 34.2, 28.6, 67.0, 0.0, 6.5, 2.0, 40, 1.0, "C"),
    @(35, "27-05-2025", 25.1, 35.5, 29.2, 69.0, 0.0, 7.5, 2.0, 300, 1.0, "C"),
    @(36, "28-05-2025", 24.7, 33.3, 28.0, 68.0, 0.0, 8.0, 3.0, 230, 1.0, "C"),
    @(37, "29-05-2025", 24.9, 33.6, 28.9, 70.0, 0.0, 6.8, 2.0, 180, 0.0, "C"),
    @(38, "30-05-2025", 25.1, 34.6, 28.6, 68.0, 0.0, 2.5, 3.0, 290, 1.0, "C"),
    @(39, "31-05-2025", 25.2, 32.7, 28.2, 71.0, 0.0, 7.8, 2.0, 40, 0.0, "C")
)

foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}

# --- 4) Leave the relocated spacer row (row 40) selected --------------------
$ws.Range("A40:K40").Select()
